$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rsquo = [char]0x2019
$a2 = "We" + $rsquo + "re Asking the Wrong Questions About AI"
$ws.Range("A2").Value = $a2

$b2 = "Okay, so this article is all about how AI is changing the game, especially for education and the workforce. Here's the gist:`n" +
      "*   AI is coming for knowledge work jobs, like coding and maybe even curriculum design!`n" +
      "*   We need to rethink what students need to know and *do* in a world where AI can do a lot of the basic cognitive stuff.`n" +
      "*   Competency-based education might be the way to go, focusing on skills and application.`n" +
      "*   There's a big opportunity to focus on uniquely human skills, like teaching, social work, and mental health support - jobs that AI can't really do.`n" +
      "*   We need to value and pay those 'human' jobs more!"
$ws.Range("B2").Value = $b2

$c2 = "Girl, this is SO relevant to our interests! Think about it: AI in education is a huge topic, and this article talks about how it's impacting curriculum and job prospects. It also brings up the importance of soft skills, which is totally your thing. Plus, the focus on valuing human-centered jobs could be a game-changer for marginalized communities if we can create pathways into those fields. And of course, it touches on how edtech and learning design need to adapt to this new AI reality. Basically, it's hitting all the key points!"
$ws.Range("C2").Value = $c2
